$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.707.40'
$ws.Range('E2').Value = '  +4.21%  '
$ws.Range('D3').Value = '2.752.19'
$ws.Range('E3').Value = '  +3.69%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '117.04'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +3.87%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '333.99'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +2.92%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.538'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +2.70%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.580'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +5.89%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '41.52'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +3.11%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '20.17'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('E13').Value = '  +2.94%  '
$ws.Range('E14').Value = '  +4.73%  '
$ws.Range('D15').Value = '3.179.04'
$ws.Range('E15').Value = '  +3.05%  '
$ws.Range('D16').Value = '2.756.52'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '51.619.27'
$ws.Range('E18').Value = '  +3.94%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '13.78'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +6.02%  '
$ws.Range('E20').Value = '  +3.31%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.87'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').Value = '  +1.40%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '278.07'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +1.92%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '70.30'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('E25').Value = '  +4.73%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '26.97'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('E29').Value = '  +2.37%  '
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '35.68'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +0.24%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '50.48'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('E34').Value = '  +3.46%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0826'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +3.28%  '
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.11'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '5.02'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '3.29'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +6.07%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '129.60'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +3.79%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '23.81'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +8.31%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.0345'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +9.68%  '
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('E45').Value = '  +3.64%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.36'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +13.61%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.40'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +4.05%  '
$ws.Range('D48').Value = '2.104.56'
$ws.Range('E48').Value = '  -0.04%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +3.52%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '5.65'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +7.07%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '9.01'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +0.55%  '
